$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '71.804.30'
Set-TextValue "E2" '  +4.52%  '

# Row 3
Set-TextValue "D3" '2.627.96'
Set-TextValue "E3" '  +4.48%  '

# Row 4
Set-TextValue "E4" '  -0.05%  '

# Row 5
Set-TextValue "D5" '607.59'
Set-TextValue "E5" '  +2.51%  '

# Row 6
Set-TextValue "D6" '179.56'
Set-TextValue "E6" '  +3.03%  '

# Row 7
Set-TextValue "E7" '  -0.08%  '

# Row 8
Set-TextValue "D8" '0.524'
Set-TextValue "E8" '  +1.64%  '

# Row 9
Set-TextValue "D9" '2.624.62'
Set-TextValue "E9" '  +4.35%  '

# Row 10
Set-TextValue "E10" '  +13.15%  '

# Row 11
Set-TextValue "E11" '  +1.22%  '

# Row 12
Set-TextValue "E12" '  +3.65%  '

# Row 13
Set-TextValue "E13" '  +0.84%  '

# Row 14
Set-TextValue "B14" 'ShibaInu'
Set-TextValue "C14" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D14" '0.0000188'
Set-TextValue "E14" '  +9.49%  '

# Row 15
Set-TextValue "B15" 'WrappedliquidstakedEther2.0'
Set-TextValue "C15" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue "D15" '3.084.17'
Set-TextValue "E15" '  +3.52%  '

# Row 16
Set-TextValue "D16" '71.751.67'
Set-TextValue "E16" '  +4.54%  '

# Row 17
Set-TextValue "D17" '26.49'
Set-TextValue "E17" '  +2.62%  '

# Row 18
Set-TextValue "D18" '2.631.20'
Set-TextValue "E18" '  +4.47%  '

# Row 19
Set-TextValue "B19" 'Uniswap'
Set-TextValue "C19" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue "D19" '8.06'
Set-TextValue "E19" '  +7.31%  '

# Row 20
Set-TextValue "B20" 'BitcoinCash'
Set-TextValue "C20" 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue "D20" '382.59'
Set-TextValue "E20" '  +5.42%  '

# Row 21
Set-TextValue "E21" '  +5.80%  '

# Row 22
Set-TextValue "D22" '4.15'
Set-TextValue "E22" '  +2.81%  '

# Row 23
Set-TextValue "B23" 'SuiNetwork'
Set-TextValue "C23" 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue "D23" '2.00'
Set-TextValue "E23" '  +21.32%  '

# Row 24
Set-TextValue "B24" 'Litecoin'
Set-TextValue "C24" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue "D24" '72.70'
Set-TextValue "E24" '  +3.42%  '

# Row 25
Set-TextValue "D25" '4.45'
Set-TextValue "E25" '  +7.18%  '

# Row 26
Set-TextValue "E26" '  +0.08%  '

# Row 27
Set-TextValue "D27" '9.97'
Set-TextValue "E27" '  +12.56%  '

# Row 28
Set-TextValue "D28" '2.758.94'
Set-TextValue "E28" '  +4.24%  '

# Row 29
Set-TextValue "D29" '0.999'
Set-TextValue "E29" '  +0.04%  '

# Row 30
Set-TextValue "D30" '0.0₃0961'
Set-TextValue "E30" '  +10.28%  '

# Row 31
Set-TextValue "D31" '547.33'
Set-TextValue "E31" '  +7.69%  '

# Row 32
Set-TextValue "D32" '8.06'
Set-TextValue "E32" '  +4.45%  '

# Row 33
Set-TextValue "E33" '  +9.08%  '

# Row 34
Set-TextValue "E34" '  +3.63%  '

# Row 35
Set-TextValue "D35" '0.999'
Set-TextValue "E35" '  -0.12%  '

# Row 36
Set-TextValue "D36" '166.35'
Set-TextValue "E36" '  +2.95%  '

# Row 37
Set-TextValue "D37" '19.21'
Set-TextValue "E37" '  +3.58%  '

# Row 38
Set-TextValue "E38" '  -2.61%  '

# Row 39
Set-TextValue "D39" '19.13'
Set-TextValue "E39" '  +2.61%  '

# Row 40
Set-TextValue "E40" '  +7.24%  '

# Row 41
Set-TextValue "D41" '1.85'
Set-TextValue "E41" '  +8.95%  '

# Row 42
Set-TextValue "B42" 'USDe'
Set-TextValue "C42" 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue "D42" '1.00'
Set-TextValue "E42" '  +0.09%  '

# Row 43
Set-TextValue "B43" 'dogwifhat'
Set-TextValue "C43" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D43" '2.62'
Set-TextValue "E43" '  +12.14%  '

# Row 44
Set-TextValue "D44" '5.03'
Set-TextValue "E44" '  +6.17%  '

# Row 45
Set-TextValue "D45" '0.332'
Set-TextValue "E45" '  +3.54%  '

# Row 46
Set-TextValue "D46" '39.38'
Set-TextValue "E46" '  +1.16%  '

# Row 47
Set-TextValue "D47" '150.78'
Set-TextValue "E47" '  +0.37%  '

# Row 48
Set-TextValue "D48" '3.64'
Set-TextValue "E48" '  +2.53%  '

# Row 49
Set-TextValue "D49" '0.539'
Set-TextValue "E49" '  +5.62%  '

# Row 50
Set-TextValue "E50" '  +8.99%  '

# Row 51
Set-TextValue "D51" '0.0₆0263'
Set-TextValue "E51" '  +5.82%  '
